$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal string (e.g. "1.00") must be written
# with a Text number format so Excel does not silently coerce them to a numeric
# value (which would also normalize away the formatting, e.g. drop "1.00" -> 1).
# The original style is captured and restored afterwards so the cell keeps no
# explicit style index, matching the source file.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '97.065.47'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.700.02'
$ws.Range('E3').Value = '  +0.19%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '2.11'
$ws.Range('E5').Value = '  +11.61%  '
Set-TextValue 'D6' '235.33'
$ws.Range('E6').Value = '  -2.01%  '
Set-TextValue 'D7' '656.88'
$ws.Range('E7').Value = '  +0.29%  '
Set-TextValue 'D8' '0.429'
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('E9').Value = '  +0.65%  '
Set-TextValue 'D10' '1.00'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.698.61'
Set-TextValue 'D12' '44.91'
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('E13').Value = '  +14.78%  '
$ws.Range('E14').Value = '  +0.38%  '
Set-TextValue 'D15' '6.84'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = '4.393.37'
$ws.Range('D17').Value = '96.772.98'
$ws.Range('E17').Value = '  +0.02%  '
Set-TextValue 'D18' '9.15'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').Value = '3.698.92'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').Value = '  +0.58%  '
Set-TextValue 'D21' '18.62'
$ws.Range('E21').Value = '  -3.41%  '
Set-TextValue 'D23' '521.57'
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('E25').Value = '  +8.59%  '
$ws.Range('E26').Value = '  -3.41%  '
Set-TextValue 'D27' '107.80'
$ws.Range('E27').Value = '  +5.07%  '
Set-TextValue 'D28' '0.201'
$ws.Range('E28').Value = '  +19.63%  '
$ws.Range('D29').Value = '3.900.10'
$ws.Range('E29').Value = '  +0.37%  '
Set-TextValue 'D30' '13.45'
$ws.Range('E30').Value = '  -0.07%  '
Set-TextValue 'D31' '12.60'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('E33').Value = '  -0.07%  '
Set-TextValue 'D34' '0.189'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('E35').Value = '  -3.49%  '
Set-TextValue 'D36' '1.00'
$ws.Range('E36').Value = '  +0.00%  '
Set-TextValue 'D37' '32.55'
$ws.Range('E37').Value = '  -0.64%  '
Set-TextValue 'D38' '638.92'
$ws.Range('E38').Value = '  -2.39%  '
Set-TextValue 'D39' '0.592'
$ws.Range('E39').Value = '  -2.64%  '
Set-TextValue 'D40' '8.69'
$ws.Range('E40').Value = '  -4.43%  '
$ws.Range('E41').Value = '  +0.02%  '
Set-TextValue 'D42' '0.166'
$ws.Range('E42').Value = '  +1.74%  '
Set-TextValue 'D43' '0.500'
$ws.Range('E43').Value = '  +9.05%  '
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('E45').Value = '  +0.83%  '
Set-TextValue 'D46' '39.79'
$ws.Range('E46').Value = '  +4.38%  '
Set-TextValue 'D47' '0.959'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('E48').Value = '  -1.02%  '
Set-TextValue 'D49' '2.38'
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('E50').Value = '  -0.16%  '
Set-TextValue 'D51' '8.73'
$ws.Range('E51').Value = '  -0.39%  '
